# Updates cryptos list values (Price column D, Volume(1h) column E)
# per commit "Updated cryptos list on Fri Aug  4 06:58:15 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.179.76'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '1.832.85'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('D4').Value = '''0.9994'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''242.08'
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('D6').Value = '''0.6581'
$ws.Range('E6').Value = '  -1.87%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '''0.07393'
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').Value = '''0.2926'
$ws.Range('E9').Value = '  -1.43%  '
$ws.Range('E10').Value = '  -0.44%  '
$ws.Range('D11').Value = '''0.07760'
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range('D12').Value = '1.835.97'
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('D13').Value = '''4.995'
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').Value = '''0.6654'
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('D15').Value = '''82.82'
$ws.Range('E15').Value = '  -3.87%  '
$ws.Range('D16').Value = '''6.101'
$ws.Range('E16').Value = '  -0.79%  '
$ws.Range('D17').Value = '''0.000008400'
$ws.Range('E17').Value = '  +1.98%  '
$ws.Range('D18').Value = '29.184.68'
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('D19').Value = '2.085.58'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').Value = '''12.44'
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('D22').Value = '''1.001'
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').Value = '''7.134'
$ws.Range('E23').Value = '  -2.48%  '
$ws.Range('D24').Value = '''1.001'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = '''158.85'
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('D26').Value = '''8.601'
$ws.Range('E26').Value = '  -0.85%  '
$ws.Range('E27').Value = '  -2.79%  '
$ws.Range('D28').Value = '''17.92'
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('E29').Value = '  +0.86%  '
$ws.Range('E30').Value = '  -2.86%  '
$ws.Range('D31').Value = '''4.042'
$ws.Range('E31').Value = '  -1.90%  '
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('D33').Value = '''0.05270'
$ws.Range('E33').Value = '  -2.33%  '
$ws.Range('D34').Value = '''1.861'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('D35').Value = '''0.7407'
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('D36').Value = '''1.139'
$ws.Range('E36').Value = '  +1.30%  '
$ws.Range('D37').Value = '''2.655'
$ws.Range('E37').Value = '  -0.95%  '
$ws.Range('D38').Value = '1.302.16'
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('D39').Value = '''0.01789'
$ws.Range('E39').Value = '  -0.85%  '
$ws.Range('D40').Value = '''2.731'
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('D41').Value = '''0.9238'
$ws.Range('E41').Value = '  -0.85%  '
$ws.Range('D42').Value = '''5.924'
$ws.Range('E42').Value = '  -2.77%  '
$ws.Range('D43').Value = '''0.08434'
$ws.Range('E43').Value = '  +2.46%  '
$ws.Range('D44').Value = '''0.9996'
$ws.Range('D45').Value = '''102.22'
$ws.Range('E45').Value = '  -2.04%  '
$ws.Range('D46').Value = '1.974.56'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('E47').Value = '  -0.72%  '
$ws.Range('E48').Value = '  -1.67%  '
$ws.Range('D49').Value = '''1.747'
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('E50').Value = '  -0.82%  '
$ws.Range('D51').Value = '''0.05847'
$ws.Range('E51').Value = '  -1.39%  '
